$wb = $excel.ActiveWorkbook

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e2cc7304e90213181a7be9dde80d917ad7fa781/e2e/c0bbf868-7e5f-4a91-8a77-517c9362d2b4.md"

# ----- Sheet "Overview" -----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "4ce20d18-9d7f-47d2-8cce-61dea857f261.md"
$ws.Range("B2").Value = "e2e\4ce20d18-9d7f-47d2-8cce-61dea857f261.md"
$ws.Range("G2").Value = "2016-08-27 12:54:32"

# Refresh hyperlink display text on B2 without altering its target relationship
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $hyperlinkAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "e2e\4ce20d18-9d7f-47d2-8cce-61dea857f261.md")

# ----- Sheet "zh-cn" -----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = "4ce20d18-9d7f-47d2-8cce-61dea857f261.md"
$ws.Range("G2").Value = "4ce20d18-9d7f-47d2-8cce-61dea857f261.5192d3e53ec7fee55a64fee9d499c3123295fba0.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-27 12:54:29"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $hyperlinkAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "4ce20d18-9d7f-47d2-8cce-61dea857f261.md")

# ----- Sheet "de-de" -----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = "4ce20d18-9d7f-47d2-8cce-61dea857f261.md"
$ws.Range("G2").Value = "4ce20d18-9d7f-47d2-8cce-61dea857f261.5192d3e53ec7fee55a64fee9d499c3123295fba0.de-de.xlf"
$ws.Range("H2").Value = "2016-08-27 12:54:32"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $hyperlinkAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "4ce20d18-9d7f-47d2-8cce-61dea857f261.md")
